$wb = $excel.ActiveWorkbook

# ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 2711.7058
$ws.Range("I82").Value = 484.53845
$ws.Range("J82").Value = 9950
$ws.Range("K82").Value = 1453.61535
$ws.Range("L82").Value = 29850
$ws.Range("M82").Value = -1047.61535
$ws.Range("N82").Value = -30662

# ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 2711.7058
$ws.Range("I85").Value = 484.53845
$ws.Range("J85").Value = 9950
$ws.Range("K85").Value = 1453.61535
$ws.Range("L85").Value = 29850
$ws.Range("M85").Value = -49.61535000000003
$ws.Range("N85").Value = -32658

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2563.5833
$ws.Range("J121").Value = 3086.4443
$ws.Range("L121").Value = 9259.332900000001
$ws.Range("N121").Value = -12753.3329

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4400.2803
$ws.Range("I132").Value = 3717.59
$ws.Range("J132").Value = 6383.3335
$ws.Range("K132").Value = 11152.77
$ws.Range("L132").Value = 19150.0005
$ws.Range("M132").Value = -8622.77
$ws.Range("N132").Value = -24210.0005

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 58029.855
$ws.Range("I136").Value = 20709
$ws.Range("J136").Value = 64250
$ws.Range("K136").Value = 20709
$ws.Range("L136").Value = 64250
$ws.Range("M136").Value = -15609
$ws.Range("N136").Value = -74450

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1380.4822
$ws.Range("I137").Value = 2522.1177
$ws.Range("J137").Value = 882.8461
$ws.Range("K137").Value = 7566.353099999999
$ws.Range("L137").Value = 2648.5383
$ws.Range("M137").Value = -5016.353099999999
$ws.Range("N137").Value = -7748.5383

# ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2038.5306
$ws.Range("I61").Value = 2278.2173
$ws.Range("J61").Value = 1826.5
$ws.Range("K61").Value = 2278.2173
$ws.Range("L61").Value = 1826.5
$ws.Range("M61").Value = -2066.2173
$ws.Range("N61").Value = -2250.5

# ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 56225
$ws.Range("J104").Value = 56225
$ws.Range("L104").Value = 56225
$ws.Range("N104").Value = -63213

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1610
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1610
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4830
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9730

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2038.5306
$ws.Range("I136").Value = 2278.2173
$ws.Range("J136").Value = 1826.5
$ws.Range("K136").Value = 6834.651899999999
$ws.Range("L136").Value = 5479.5
$ws.Range("M136").Value = -4284.651899999999
$ws.Range("N136").Value = -10579.5

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 28098.076
$ws.Range("I139").Value = 10000
$ws.Range("J139").Value = 29606.25
$ws.Range("K139").Value = 10000
$ws.Range("L139").Value = 29606.25
$ws.Range("M139").Value = -4860
$ws.Range("N139").Value = -39886.25

# BSM row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 23641.334
$ws.Range("J95").Value = 23641.334
$ws.Range("L95").Value = 23641.334
$ws.Range("N95").Value = -29133.334

# CRP row 43
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 19325.6
$ws.Range("J43").Value = 19325.6
$ws.Range("L43").Value = 19325.6
$ws.Range("N43").Value = -19693.6

# CRP row 96
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 27475
$ws.Range("J96").Value = 27475
$ws.Range("L96").Value = 27475
$ws.Range("N96").Value = -32967

# CRP row 101
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 19325.6
$ws.Range("J101").Value = 19325.6
$ws.Range("L101").Value = 19325.6
$ws.Range("N101").Value = -25815.6

# CRP row 111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 49800
$ws.Range("J111").Value = 49800
$ws.Range("L111").Value = 49800
$ws.Range("N111").Value = -57980

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 100001070
$ws.Range("I122").Value = 166667410
$ws.Range("J122").Value = 1578.5
$ws.Range("K122").Value = 500002230
$ws.Range("L122").Value = 4735.5
$ws.Range("M122").Value = -499999780
$ws.Range("N122").Value = -9635.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2866.1592
$ws.Range("I132").Value = 2451
$ws.Range("J132").Value = 3103.3928
$ws.Range("K132").Value = 7353
$ws.Range("L132").Value = 9310.178400000001
$ws.Range("M132").Value = -4823
$ws.Range("N132").Value = -14370.1784

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1334.7142
$ws.Range("I81").Value = 914
$ws.Range("J81").Value = 1355.75
$ws.Range("K81").Value = 2742
$ws.Range("L81").Value = 4067.25
$ws.Range("M81").Value = -1619
$ws.Range("N81").Value = -6313.25

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1334.7142
$ws.Range("I84").Value = 914
$ws.Range("J84").Value = 1355.75
$ws.Range("K84").Value = 8226
$ws.Range("L84").Value = 12201.75
$ws.Range("M84").Value = -2610
$ws.Range("N84").Value = -23433.75

# CUL row 105
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 10389
$ws.Range("J105").Value = 10984.889
$ws.Range("L105").Value = 32954.667
$ws.Range("N105").Value = -38196.667

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2556.7144
$ws.Range("I132").Value = 1979.64
$ws.Range("J132").Value = 3157.8333
$ws.Range("K132").Value = 5938.92
$ws.Range("L132").Value = 9473.499899999999
$ws.Range("M132").Value = -3408.92
$ws.Range("N132").Value = -14533.4999

# GSM row 141
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 46038.168
$ws.Range("J141").Value = 46038.168
$ws.Range("L141").Value = 46038.168
$ws.Range("N141").Value = -56398.168

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 22633.166
$ws.Range("I122").Value = 38666.668
$ws.Range("J122").Value = 6599.6665
$ws.Range("K122").Value = 116000.004
$ws.Range("L122").Value = 19798.9995
$ws.Range("M122").Value = -113550.004
$ws.Range("N122").Value = -24698.9995

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2524.818
$ws.Range("I122").Value = 1984.4667
$ws.Range("J122").Value = 3682.7144
$ws.Range("K122").Value = 5953.4001
$ws.Range("L122").Value = 11048.1432
$ws.Range("M122").Value = -3503.4001
$ws.Range("N122").Value = -15948.1432

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1810.8334
$ws.Range("I132").Value = 1686.4333
$ws.Range("J132").Value = 1966.3334
$ws.Range("K132").Value = 5059.2999
$ws.Range("L132").Value = 5899.0002
$ws.Range("M132").Value = -2529.2999
$ws.Range("N132").Value = -10959.0002

# WVR row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 49460.555
$ws.Range("J141").Value = 49460.555
$ws.Range("L141").Value = 49460.555
$ws.Range("N141").Value = -59820.555

Write-Output "Applied all Mandragora_Profits updates."
